$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Data")

# New row 4 data (fraud_orders first, establishes shared-string order)
$ws.Range("A4").Value = "fraud_orders"

# New header cells
$ws.Range("T1").Value = "FraudCancelOrderConfirmationId"
$ws.Range("U1").Value = "FraudReprocessOrderConfirmationId"

$ws.Range("T4").Value = "US-C-23080455071226315254"
$ws.Range("U4").Value = "US-C-23080226523476649177"
